$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5ec0c1eb1c008354f9286e488c21f32ca8fad9b/e2e/4b4c19a2-07bd-4260-a3ed-2607b79109fd.md"
$mdFileName = "4b4c19a2-07bd-4260-a3ed-2607b79109fd.md"
$newStatusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status mirrors for zh-cn / de-de both flip to "handed back"
# ---------------------------------------------------------------------------
$ws1.Range("E2").Value = $newStatusText
$ws1.Range("F2").Value = $newStatusText

# Overview columns E (zh-cn) and F (de-de) grow wider to fit the longer text
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet (report for handback)
# ---------------------------------------------------------------------------
$ws2.Range("C2").Value = $newStatusText

$ws2.Range("I2").Value = $mdFileName
$ws2.Hyperlinks.Add($ws2.Range("I2"), $hyperlinkTarget, "", "", $mdFileName) | Out-Null
$ws2.Range("I2").Font.Underline = $True
$ws2.Range("I2").Font.Color = 15570276

$ws2.Range("J2").Value = "4b4c19a2-07bd-4260-a3ed-2607b79109fd.1abd6df45ed39ce7c1f45640f125c11961bf1118.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-02 01:10:20"

$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws2.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws2.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet (report for handback)
# ---------------------------------------------------------------------------
$ws3.Range("C2").Value = $newStatusText

$ws3.Range("I2").Value = $mdFileName
$ws3.Hyperlinks.Add($ws3.Range("I2"), $hyperlinkTarget, "", "", $mdFileName) | Out-Null
$ws3.Range("I2").Font.Underline = $True
$ws3.Range("I2").Font.Color = 15570276

$ws3.Range("J2").Value = "4b4c19a2-07bd-4260-a3ed-2607b79109fd.1abd6df45ed39ce7c1f45640f125c11961bf1118.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-02 01:10:27"

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws3.Columns.Item(9).ColumnWidth = 39.166666666666664
$ws3.Columns.Item(10).ColumnWidth = 39.166666666666664
